# Update "想去人数" (want-to-go count) figures on the two sheets that carry
# this exhibition listing data: "展览" and "全部类型" (identical content).
#
# Row -> (column F old -> new)
#   3  : 3119 -> 3121
#   17 :  230 ->  231
#   22 :   17 ->   18
#   28 :   73 ->   75
#   29 : 2094 -> 2095
#   32 :  461 ->  462
#   33 :  199 ->  201

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3121
    17 = 231
    22 = 18
    28 = 75
    29 = 2095
    32 = 462
    33 = 201
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
